# Regenerate orders with updated distance/size codes.
#
# The Distance codes (D64 -> D69, D80 -> D86, D51 -> D55) and the S30 size
# code (-> S31) changed throughout the trial table. Every other string
# (Face ids, S25/S20, fixation/condition labels, etc.) stays the same.
# We apply the substitutions as plain text replacements across the used
# range so any cell built from these tokens (Condition, Filename_Left,
# Filename_Right, Distance, Size) gets updated consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$firstRow = $used.Row
$firstCol = $used.Column

# Order matters: longer/more-specific tokens first so substrings don't
# collide (none of these actually overlap, but keep it safe & explicit).
$replacements = [ordered]@{
    'D64' = 'D69'
    'D80' = 'D86'
    'D51' = 'D55'
    'S30' = 'S31'
}

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($firstRow + $r, $firstCol + $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            $newVal = $val
            foreach ($key in $replacements.Keys) {
                $newVal = $newVal.Replace($key, $replacements[$key])
            }
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
